# Apply cell value updates described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'68.094.85"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.64%  "
$c = $ws.Range("D3")
$c.Value = "'3.790.57"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.20%  "
$c = $ws.Range("D4")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$c = $ws.Range("D5")
$c.Value = "'599.86"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "
$c = $ws.Range("D6")
$c.Value = "'164.31"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.62%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("E10").Value = "  -0.02%  "
$c = $ws.Range("D11")
$c.Value = "'6.55"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +3.91%  "
$ws.Range("E12").Value = "  -2.11%  "
$c = $ws.Range("D13")
$c.Value = "'35.59"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.84%  "
$c = $ws.Range("D14")
$c.Value = "'4.425.59"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.39%  "
$c = $ws.Range("D15")
$c.Value = "'3.822.68"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.03%  "
$c = $ws.Range("D16")
$c.Value = "'68.075.48"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.63%  "
$c = $ws.Range("D17")
$c.Value = "'18.29"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.40%  "
$ws.Range("E18").Value = "  +2.29%  "
$ws.Range("E19").Value = "  -0.46%  "
$c = $ws.Range("D20")
$c.Value = "'460.83"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.04%  "
$c = $ws.Range("D21")
$c.Value = "'9.65"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.58%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("E23").Value = "  -4.24%  "
$c = $ws.Range("D24")
$c.Value = "'82.91"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.61%  "
$c = $ws.Range("D25")
$c.Value = "'11.97"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.02%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -0.04%  "
$c = $ws.Range("D28")
$c.Value = "'9.98"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.19%  "
$c = $ws.Range("D29")
$c.Value = "'3.939.48"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  -5.11%  "
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("E35").Value = "  -0.84%  "
$ws.Range("E36").Value = "  -0.14%  "
$c = $ws.Range("D37")
$c.Value = "'0.140"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.35%  "
$c = $ws.Range("D38")
$c.Value = "'3.29"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.46%  "
$ws.Range("E39").Value = "  +0.72%  "
$c = $ws.Range("D40")
$c.Value = "'0.984"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.36%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E43").Value = "  -1.55%  "
$c = $ws.Range("D44")
$c.Value = "'43.45"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("E45").Value = "  -0.62%  "
$c = $ws.Range("D46")
$c.Value = "'152.10"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.81%  "
$c = $ws.Range("D47")
$c.Value = "'8.33"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("E48").Value = "  +1.47%  "
$ws.Range("E49").Value = "  +0.34%  "
$c = $ws.Range("D50")
$c.Value = "'389.22"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.32%  "
$c = $ws.Range("D51")
$c.Value = "'26.48"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.35%  "
